$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4018.6667
$ws.Range("I62").Value = 4223.6665
$ws.Range("J62").Value = 3198.6667
$ws.Range("K62").Value = 4223.6665
$ws.Range("L62").Value = 3198.6667
$ws.Range("M62").Value = -3599.6665
$ws.Range("N62").Value = -4446.6667
$ws.Range("H65").Value = 4018.6667
$ws.Range("I65").Value = 4223.6665
$ws.Range("J65").Value = 3198.6667
$ws.Range("K65").Value = 21118.3325
$ws.Range("L65").Value = 15993.3335
$ws.Range("M65").Value = -17998.3325
$ws.Range("N65").Value = -22233.3335
$ws.Range("H125").Value = 77992.92999999999
$ws.Range("I125").Value = 168800.17
$ws.Range("J125").Value = 9887.5
$ws.Range("K125").Value = 1519201.53
$ws.Range("L125").Value = 88987.5
$ws.Range("M125").Value = -1516741.53
$ws.Range("N125").Value = -93907.5
$ws.Range("H137").Value = 9092644
$ws.Range("I137").Value = 15153410
$ws.Range("J137").Value = 1495.5454
$ws.Range("K137").Value = 45460230
$ws.Range("L137").Value = 4486.6362
$ws.Range("M137").Value = -45457680
$ws.Range("N137").Value = -9586.636200000001
$ws.Range("H138").Value = 2208.132
$ws.Range("J138").Value = 2425.1292
$ws.Range("L138").Value = 7275.3876
$ws.Range("N138").Value = -17555.3876

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5140773.5
$ws.Range("I32").Value = 7226.7964
$ws.Range("J32").Value = 30341822
$ws.Range("K32").Value = 7226.7964
$ws.Range("L32").Value = 30341822
$ws.Range("M32").Value = -6939.7964
$ws.Range("N32").Value = -30342396
$ws.Range("H61").Value = 2250.45
$ws.Range("I61").Value = 1447.2778
$ws.Range("J61").Value = 2907.5908
$ws.Range("K61").Value = 1447.2778
$ws.Range("L61").Value = 2907.5908
$ws.Range("M61").Value = -1235.2778
$ws.Range("N61").Value = -3331.5908
$ws.Range("H74").Value = 25424520
$ws.Range("I74").Value = 33334126
$ws.Range("J74").Value = 789
$ws.Range("K74").Value = 33334126
$ws.Range("L74").Value = 789
$ws.Range("M74").Value = -33333252
$ws.Range("N74").Value = -2537
$ws.Range("H77").Value = 25424520
$ws.Range("I77").Value = 33334126
$ws.Range("J77").Value = 789
$ws.Range("K77").Value = 166670630
$ws.Range("L77").Value = 3945
$ws.Range("M77").Value = -166666262
$ws.Range("N77").Value = -12681
$ws.Range("H88").Value = 1666.6666
$ws.Range("I88").Value = 1666.6666
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 1666.6666
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -1260.6666
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 1666.6666
$ws.Range("I91").Value = 1666.6666
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 1666.6666
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -262.6666
$ws.Range("N91").ClearContents()
$ws.Range("H132").Value = 3129805
$ws.Range("I132").Value = 2558.625
$ws.Range("J132").Value = 5214636
$ws.Range("K132").Value = 7675.875
$ws.Range("L132").Value = 15643908
$ws.Range("M132").Value = -5145.875
$ws.Range("N132").Value = -15648968
$ws.Range("H136").Value = 2250.45
$ws.Range("I136").Value = 1447.2778
$ws.Range("J136").Value = 2907.5908
$ws.Range("K136").Value = 4341.8334
$ws.Range("L136").Value = 8722.7724
$ws.Range("M136").Value = -1791.8334
$ws.Range("N136").Value = -13822.7724

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1613.9333
$ws.Range("I86").Value = 1384.4546
$ws.Range("J86").Value = 2245
$ws.Range("K86").Value = 1384.4546
$ws.Range("L86").Value = 2245
$ws.Range("M86").Value = -261.4546
$ws.Range("N86").Value = -4491
$ws.Range("H89").Value = 1613.9333
$ws.Range("I89").Value = 1384.4546
$ws.Range("J89").Value = 2245
$ws.Range("K89").Value = 6922.273
$ws.Range("L89").Value = 11225
$ws.Range("M89").Value = -1306.273
$ws.Range("N89").Value = -22457

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6495315
$ws.Range("I31").Value = 1435.2106
$ws.Range("J31").Value = 12822685
$ws.Range("K31").Value = 1435.2106
$ws.Range("L31").Value = 12822685
$ws.Range("M31").Value = -1140.2106
$ws.Range("N31").Value = -12823275
$ws.Range("H34").Value = 6495315
$ws.Range("I34").Value = 1435.2106
$ws.Range("J34").Value = 12822685
$ws.Range("K34").Value = 1435.2106
$ws.Range("L34").Value = 12822685
$ws.Range("M34").Value = -1233.2106
$ws.Range("N34").Value = -12823089
$ws.Range("H107").Value = 1656.05
$ws.Range("I107").Value = 671.46155
$ws.Range("J107").Value = 3484.5715
$ws.Range("K107").Value = 671.46155
$ws.Range("L107").Value = 3484.5715
$ws.Range("M107").Value = 1248.53845
$ws.Range("N107").Value = -7324.5715
$ws.Range("H132").Value = 3036.1667
$ws.Range("J132").Value = 3374.25
$ws.Range("L132").Value = 10122.75
$ws.Range("N132").Value = -15182.75
$ws.Range("H134").Value = 2785.875
$ws.Range("I134").Value = 1182.4375
$ws.Range("J134").Value = 5992.75
$ws.Range("K134").Value = 3547.3125
$ws.Range("L134").Value = 17978.25
$ws.Range("M134").Value = -1012.3125
$ws.Range("N134").Value = -23048.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 11357.889
$ws.Range("I23").Value = 95.5
$ws.Range("J23").Value = 14575.714
$ws.Range("K23").Value = 286.5
$ws.Range("L23").Value = 43727.142
$ws.Range("M23").Value = -51.5
$ws.Range("N23").Value = -44197.142
$ws.Range("H68").Value = 914.96
$ws.Range("I68").Value = 808.2045000000001
$ws.Range("J68").Value = 998.8393
$ws.Range("K68").Value = 2424.6135
$ws.Range("L68").Value = 2996.5179
$ws.Range("M68").Value = -1613.6135
$ws.Range("N68").Value = -4618.5179
$ws.Range("H71").Value = 914.96
$ws.Range("I71").Value = 808.2045000000001
$ws.Range("J71").Value = 998.8393
$ws.Range("K71").Value = 7273.8405
$ws.Range("L71").Value = 8989.5537
$ws.Range("M71").Value = -3217.8405
$ws.Range("N71").Value = -17101.5537
$ws.Range("H107").Value = 15471742
$ws.Range("I107").Value = 29413378
$ws.Range("J107").Value = 658753.7
$ws.Range("K107").Value = 88240134
$ws.Range("L107").Value = 1976261.1
$ws.Range("M107").Value = -88238214
$ws.Range("N107").Value = -1980101.1
$ws.Range("H113").Value = 473.14285
$ws.Range("I113").Value = 424.375
$ws.Range("J113").Value = 538.1667
$ws.Range("K113").Value = 1273.125
$ws.Range("L113").Value = 1614.5001
$ws.Range("M113").Value = 896.875
$ws.Range("N113").Value = -5954.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6576.923
$ws.Range("I70").Value = 6666.6665
$ws.Range("J70").Value = 6550
$ws.Range("K70").Value = 6666.6665
$ws.Range("L70").Value = 6550
$ws.Range("M70").Value = -6396.6665
$ws.Range("N70").Value = -7090
$ws.Range("H73").Value = 6576.923
$ws.Range("I73").Value = 6666.6665
$ws.Range("J73").Value = 6550
$ws.Range("K73").Value = 6666.6665
$ws.Range("L73").Value = 6550
$ws.Range("M73").Value = -5730.6665
$ws.Range("N73").Value = -8422
$ws.Range("H102").Value = 2648248.5
$ws.Range("I102").Value = 4466242
$ws.Range("J102").Value = 3894.182
$ws.Range("K102").Value = 4466242
$ws.Range("L102").Value = 3894.182
$ws.Range("M102").Value = -4464620
$ws.Range("N102").Value = -7138.182
$ws.Range("H132").Value = 3227
$ws.Range("I132").Value = 2957.889
$ws.Range("J132").Value = 4034.3333
$ws.Range("K132").Value = 8873.667000000001
$ws.Range("L132").Value = 12102.9999
$ws.Range("M132").Value = -6343.667000000001
$ws.Range("N132").Value = -17162.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2066.6924
$ws.Range("I7").Value = 2270.5715
$ws.Range("J7").Value = 1828.8334
$ws.Range("K7").Value = 2270.5715
$ws.Range("L7").Value = 1828.8334
$ws.Range("M7").Value = -2158.5715
$ws.Range("N7").Value = -2052.8334
$ws.Range("H126").Value = 2066.6924
$ws.Range("I126").Value = 2270.5715
$ws.Range("J126").Value = 1828.8334
$ws.Range("K126").Value = 6811.7145
$ws.Range("L126").Value = 5486.5002
$ws.Range("M126").Value = -4341.7145
$ws.Range("N126").Value = -10426.5002
$ws.Range("H136").Value = 15626181
$ws.Range("I136").Value = 26316624
$ws.Range("J136").Value = 1688.7693
$ws.Range("K136").Value = 78949872
$ws.Range("L136").Value = 5066.3079
$ws.Range("M136").Value = -78947322
$ws.Range("N136").Value = -10166.3079

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 8086.5835
$ws.Range("I122").Value = 9424.352999999999
$ws.Range("J122").Value = 4837.7144
$ws.Range("K122").Value = 28273.059
$ws.Range("L122").Value = 14513.1432
$ws.Range("M122").Value = -25823.059
$ws.Range("N122").Value = -19413.1432
$ws.Range("H136").Value = 8622115
$ws.Range("I136").Value = 10417689
$ws.Range("J136").Value = 3364
$ws.Range("K136").Value = 31253067
$ws.Range("L136").Value = 10092
$ws.Range("M136").Value = -31250517
$ws.Range("N136").Value = -15192

